# Regenerate save_data: replace Strike# values in column G ("K") with the
# recomputed strikeout counts (s_vals) for the brogdon_connor workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new value for column G ("K")
$updates = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 0
    6  = 2
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 2
    12 = 0
    13 = 0
    14 = 0
    15 = 3
    16 = 1
    17 = 1
    18 = 0
    19 = 2
    20 = 2
    21 = 0
    22 = 0
    23 = 1
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 1
    29 = 2
    30 = 2
    31 = 1
    33 = 2
    34 = 1
    35 = 2
    36 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
